$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row = 88
        A = "Globo"
        B = "RJ TV 2"
        C = "Infraestrutura"
        D = "2025-04-08T19:16"
        E = "Neutro"
        F = "Média de 14 mil furtos de bicicletas por dia no RJ. Instituto de Segurança Pública registra aumento de 59% nos dois primeiros meses do ano. Repórter *ao vivo* do Rio de Janeiro. Em Campos, há uma ciclovia que corta a cidade indo do Parque Imperial até perto da BR-101. Em 2024, região norte fluminense registrou 38 furtos, sendo 35 em Campos. Este ano, 36 na região e 29 em Campos. No Estado, foram 526 furtos ano passado. Este ano, 848 no mesmo período deste ano."
    },
    @{
        Row = 89
        A = "Globo"
        B = "RJ TV 2"
        C = "Infraestrutura"
        D = "2025-04-08T19:30"
        E = "Negativo"
        F = "Por whatsApp, moradora de Campos solicita ajuda à emissora para cobrar do poder público a limpeza de cabos e fios pendurados, atrapalhando mobilidade de idosos e gestantes, principalmente, e todos os outros pedestres. Foram exibidas imagens da Rua Manoel Teodoro, na Pelinca, e na Rua Gil de Góis próximo ao Centro de Saúde. Apresentador lembrou que já teve uma ação em Campos para retirada destes fios.  "
    },
    @{
        Row = 90
        A = "Globo"
        B = "RJ TV 2"
        C = "CCZ"
        D = "2025-04-08T19:31"
        E = "Neutro"
        F = "Por whatsApp, chegou mensagem mostrando a piscina de uma casa que fica entre as ruas Voluntários da Pátria e Dr. Siqueira. Ela questiona como querem que acabe a epidemia do aedes aegypti. Não houve reclamação direta à prefeitura. "
    }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
